$d = $word.ActiveDocument

# --- First occurrence: "For i = 2 To 70926" -> "For i = 2 To 760192" (stockdata_easy) ---
$rng1 = $d.Content
$rng1.Find.Execute("70926", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Delete()
$ins1 = $d.Range($rng1.Start, $rng1.Start)
$ins1.InsertAfter("760192")

# --- Second occurrence: "For i = 2 To 70926" -> "For i = 2 To 760192" (stockdata_moderate) ---
$rng2 = $d.Content
$rng2.Find.Execute("70926", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Delete()
$ins2 = $d.Range($rng2.Start, $rng2.Start)
$ins2.InsertAfter("760192")

# This second edit is the most recent change, so Word's "_GoBack" bookmark should now
# mark this location (right after the freshly typed "760192") instead of staying at
# the end of the document where it originally lived.
$endPos = $ins2.Start + 6   # length of "760192"

# A bookmark placed exactly at a paragraph's final text offset needs a non-boundary
# position to anchor against, so insert a temporary marker character right after the
# edit, carve out a true zero-width point range there, then remove the marker again.
$tmp = $d.Range($endPos, $endPos)
$tmp.InsertAfter("X")

$goBackPoint = $d.Range($endPos, $endPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $goBackPoint)

$d.Range($endPos, $endPos + 1).Delete()
